# 042922_titer_calculations.xlsx — reorder Library_A/Library_B rows so B
# precedes A within each infection-sample pair, and freeze the Titers
# (TU/mL) column to its last-calculated values instead of a live formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap each Library_A / Library_B row pair (columns A:E) -------------
# Row 2 <-> Row 3 : Library_A_HPs_concentrated / Library_B_HPs_concentrated
$r1 = $ws.Range("A2:E2")
$r2 = $ws.Range("A3:E3")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

# Row 4 <-> Row 5 : Library_A_HPs_unconcentrated / Library_B_HPs_unconcentrated
$r1 = $ws.Range("A4:E4")
$r2 = $ws.Range("A5:E5")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

# Row 6 <-> Row 7 : Library_A_HPs+VSVG / Library_B_HPs+VSVG
$r1 = $ws.Range("A6:E6")
$r2 = $ws.Range("A7:E7")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

# Row 8 <-> Row 9 : Library_A_No_Plasmids / Library_B_No_Plasmids
$r1 = $ws.Range("A8:E8")
$r2 = $ws.Range("A9:E9")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

# --- Column E (Titers_(TU/mL)) is now plain numbers, not live formulas --
# The row-swap above already replaced the formulas with their computed
# results (Value, not Formula, was copied), so E2:E9 are static now.
# Drop the now-unused explicit 2-decimal number format that only existed
# to dress up the formula cells.
$ws.Range("E2:E9").Style = "Normal"

# --- Restore the on-save cursor position --------------------------------
$null = $ws.Range("B13").Select()
